# Fixed (mostly) aspect ratio issues
# Update the bounding box coordinates on the "Mapping" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -75.8787
$ws.Range("B2").Value = -75.8058

$ws.Range("A3").Value = 38.9503
$ws.Range("B3").Value = 39.0063

$ws.Range("A4").Value = -75.1437
$ws.Range("B4").Value = -75.2173

$ws.Range("A5").Value = 39.5112
$ws.Range("B5").Value = 39.4554
